# Weekly refresh of the Espinaca price series: the per-row market figures
# (Fecha, Volumen, Precio minimo/maximo/promedio, Origen, Precio $/Kg) get
# reshuffled across the existing data rows (2-24). Columns that identify the
# market/product (A,B,C,E,F,G,H,I,N,Q,R) stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as one "row record".
$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot every source row's record BEFORE any write, so the permutation
# below is safe regardless of write order (rows 9, 13, 14 map to themselves).
$snapshot = @{}
for ($r = 2; $r -le 24; $r++) {
    $rec = @{}
    foreach ($c in $cols) {
        $rec[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rec
}

# Destination row -> source row (data that now lives in that row).
$mapping = @{
    2  = 21
    3  = 12
    4  = 8
    5  = 23
    6  = 20
    7  = 10
    8  = 5
    9  = 9
    10 = 2
    11 = 4
    12 = 16
    13 = 13
    14 = 14
    15 = 24
    16 = 18
    17 = 15
    18 = 22
    19 = 11
    20 = 19
    21 = 3
    22 = 6
    23 = 7
    24 = 17
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rec = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $rec[$c]
    }
}
